$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 50; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $old = $cell.Value2
    $new = [Math]::Min($old - 1, 4)
    $cell.Value2 = $new
}

$ws.Range("H16").Select()
